$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1166.6666
$ws.Range("I58").Value = 500
$ws.Range("K58").Value = 1500
$ws.Range("M58").Value = -1350
$ws.Range("H62").Value = 4249.5
$ws.Range("J62").Value = 4249.5
$ws.Range("L62").Value = 4249.5
$ws.Range("N62").Value = -5497.5
$ws.Range("H65").Value = 4249.5
$ws.Range("J65").Value = 4249.5
$ws.Range("L65").Value = 21247.5
$ws.Range("N65").Value = -27487.5
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("H86").Value = 2339.8
$ws.Range("I86").Value = 2050
$ws.Range("K86").Value = 2050
$ws.Range("M86").Value = -927
$ws.Range("H89").Value = 2339.8
$ws.Range("I89").Value = 2050
$ws.Range("K89").Value = 10250
$ws.Range("M89").Value = -4634
$ws.Range("H96").Value = 18404.334
$ws.Range("I96").Value = 26106.5
$ws.Range("K96").Value = 78319.5
$ws.Range("M96").Value = -76946.5
$ws.Range("H111").Value = 7107.7646
$ws.Range("I111").Value = 7571.5835
$ws.Range("K111").Value = 22714.7505
$ws.Range("M111").Value = -19647.7505
$ws.Range("H131").Value = 31074.75
$ws.Range("I131").Value = 37398
$ws.Range("K131").Value = 112194
$ws.Range("M131").Value = -107154
$ws.Range("H135").Value = 4776.8335
$ws.Range("I135").Value = 4065.125
$ws.Range("J135").Value = 6200.25
$ws.Range("K135").Value = 36586.125
$ws.Range("L135").Value = 55802.25
$ws.Range("M135").Value = -34051.125
$ws.Range("N135").Value = -60872.25

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 4820
$ws.Range("J46").Value = 4820
$ws.Range("L46").Value = 4820
$ws.Range("N46").Value = -5458
$ws.Range("H92").Value = 19750
$ws.Range("J92").Value = 19750
$ws.Range("L92").Value = 19750
$ws.Range("N92").Value = -24742
$ws.Range("H97").Value = 954.4545000000001
$ws.Range("I97").Value = 901.4
$ws.Range("J97").Value = 1485
$ws.Range("K97").Value = 901.4
$ws.Range("L97").Value = 1485
$ws.Range("M97").Value = -405.4
$ws.Range("N97").Value = -2477
$ws.Range("H102").Value = 900
$ws.Range("I102").Value = 900
$ws.Range("K102").Value = 900
$ws.Range("M102").Value = 722
$ws.Range("H122").Value = 1075

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = ""

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2383.4736
$ws.Range("I16").Value = 2463.6667
$ws.Range("J16").Value = 2311.3
$ws.Range("K16").Value = 2463.6667
$ws.Range("L16").Value = 2311.3
$ws.Range("M16").Value = -2176.6667
$ws.Range("N16").Value = -2885.3
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = ""
$ws.Range("H99").Value = 2068.125
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2068.125
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 2068.125
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = -5064.125
$ws.Range("H105").Value = 3456.25
$ws.Range("I105").Value = 2392.8
$ws.Range("K105").Value = 2392.8
$ws.Range("M105").Value = -645.8000000000002
$ws.Range("H113").Value = 2383.4736
$ws.Range("I113").Value = 2463.6667
$ws.Range("J113").Value = 2311.3
$ws.Range("K113").Value = 2463.6667
$ws.Range("L113").Value = 2311.3
$ws.Range("M113").Value = -293.6667000000002
$ws.Range("N113").Value = -6651.3
$ws.Range("H122").Value = 1732.6666
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 1479.2
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 4437.6
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -9337.6
$ws.Range("H126").Value = 2068.125
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2068.125
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6204.375
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -11144.375

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 38461
$ws.Range("I63").Value = 38461
$ws.Range("K63").Value = 115383
$ws.Range("M63").Value = -114634
$ws.Range("H66").Value = 38461
$ws.Range("I66").Value = 38461
$ws.Range("K66").Value = 346149
$ws.Range("M66").Value = -342405
$ws.Range("H95").Value = 8892.5
$ws.Range("J95").Value = 8892.5
$ws.Range("L95").Value = 26677.5
$ws.Range("N95").Value = -30795.5
$ws.Range("H114").Value = 1985.6
$ws.Range("I114").Value = 28
$ws.Range("J114").Value = 2475
$ws.Range("K114").Value = 84
$ws.Range("L114").Value = 7425
$ws.Range("M114").Value = 3170
$ws.Range("N114").Value = -13933
$ws.Range("H120").Value = 799.5
$ws.Range("I120").Value = 799.5
$ws.Range("K120").Value = 2398.5
$ws.Range("M120").Value = 2439.5

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2164
$ws.Range("I102").Value = 2996
$ws.Range("J102").Value = 1997.6
$ws.Range("K102").Value = 2996
$ws.Range("L102").Value = 1997.6
$ws.Range("M102").Value = -1374
$ws.Range("N102").Value = -5241.6
$ws.Range("H122").Value = 3260.7334
$ws.Range("I122").Value = 3377.7693
$ws.Range("K122").Value = 10133.3079
$ws.Range("M122").Value = -7683.3079
$ws.Range("H132").Value = 6518.5356
$ws.Range("I132").Value = 6849.087
$ws.Range("K132").Value = 20547.261
$ws.Range("M132").Value = -18017.261

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10137.77
$ws.Range("I7").Value = 13333
$ws.Range("J7").Value = 9179.200000000001
$ws.Range("K7").Value = 13333
$ws.Range("L7").Value = 9179.200000000001
$ws.Range("M7").Value = -13221
$ws.Range("N7").Value = -9403.200000000001
$ws.Range("H22").Value = 5500
$ws.Range("I22").Value = 5500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 5500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -5205
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 5500
$ws.Range("I27").Value = 5500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 5500
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -5393
$ws.Range("N27").Value = ""
$ws.Range("H40").Value = 5750
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H122").Value = 10093.529
$ws.Range("I122").Value = 14400.8
$ws.Range("J122").Value = 8298.833000000001
$ws.Range("K122").Value = 43202.39999999999
$ws.Range("L122").Value = 24896.499
$ws.Range("M122").Value = -40752.39999999999
$ws.Range("N122").Value = -29796.499
$ws.Range("H126").Value = 10137.77
$ws.Range("I126").Value = 13333
$ws.Range("J126").Value = 9179.200000000001
$ws.Range("K126").Value = 39999
$ws.Range("L126").Value = 27537.6
$ws.Range("M126").Value = -37529
$ws.Range("N126").Value = -32477.6

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = ""
$ws.Range("H122").Value = 2500
$ws.Range("H126").Value = 2592.2
$ws.Range("I126").Value = 2003.1428
$ws.Range("K126").Value = 6009.428400000001
$ws.Range("M126").Value = -3539.428400000001
